# "Final Changes to Data and README"
#
# The diff swaps the paired "Bench" / "Squat" columns in the two header
# rows (row 2 and row 38) of Sheet1: column A now shows what used to be
# in column B (and vice versa), same for C/D and for the mirrored
# G/H/I/J block. It also moves the saved selection/scroll position and
# nudges the best-fit widths of columns A/B/G/H to follow the swapped
# (now wider/narrower) header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 header swap -------------------------------------------------
# A2 <-> B2 ("Number of Benches" <-> "Number of Squat Racks")
$ws.Range("A2").Value = "Number of Squat Racks"
$ws.Range("B2").Value = "Number of Benches"
# C2 <-> D2 ("Average Bench Line Length" <-> "Average Squat Line Length")
$ws.Range("C2").Value = "Average Squat Line Length"
$ws.Range("D2").Value = "Average Bench Line Length"
# Mirrored block: G2 <-> H2, I2 <-> J2
$ws.Range("G2").Value = "Number of Squat Racks"
$ws.Range("H2").Value = "Number of Benches"
$ws.Range("I2").Value = "Average Squat Line Length"
$ws.Range("J2").Value = "Average Bench Line Length"

# --- Row 38 header swap (identical pattern) ----------------------------
$ws.Range("A38").Value = "Number of Squat Racks"
$ws.Range("B38").Value = "Number of Benches"
$ws.Range("C38").Value = "Average Squat Line Length"
$ws.Range("D38").Value = "Average Bench Line Length"
$ws.Range("G38").Value = "Number of Squat Racks"
$ws.Range("H38").Value = "Number of Benches"
$ws.Range("I38").Value = "Average Squat Line Length"
$ws.Range("J38").Value = "Average Bench Line Length"

# --- Column widths follow the swapped header text (best-fit) ----------
# Column A/G now hold the longer "Number of Squat Racks" text (was in B/H)
$ws.Columns.Item(1).ColumnWidth = 20.833333333333332
$ws.Columns.Item(7).ColumnWidth = 20.833333333333332
# Column H now holds the shorter "Number of Benches" text (was in A)
$ws.Columns.Item(8).ColumnWidth = 17.833333333333332
# Column B now holds "Number of Benches" too, re-fitted slightly narrower
$ws.Columns.Item(2).ColumnWidth = 18.333333333333332

# --- Selection / scroll position / zoom --------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 76
$win.ScrollColumn = 1
$ws.Range("G93").Select()
$win.Zoom = 100
